$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 2370.0908
$ws.Range("J97").Value = 3463.4285
$ws.Range("L97").Value = 10390.2855
$ws.Range("N97").Value = -11382.2855
$ws.Range("H113").Value = 5925.619
$ws.Range("I113").Value = 7154.7
$ws.Range("K113").Value = 7154.7
$ws.Range("M113").Value = -3900.7
$ws.Range("H125").Value = 5075.9
$ws.Range("I125").Value = 4079.4
$ws.Range("J125").Value = 6072.4
$ws.Range("K125").Value = 36714.6
$ws.Range("L125").Value = 54651.6
$ws.Range("M125").Value = -34254.6
$ws.Range("N125").Value = -59571.6
$ws.Range("H127").Value = 1931.1
$ws.Range("I127").Value = 2274.1428
$ws.Range("J127").Value = 1746.3846
$ws.Range("K127").Value = 6822.428400000001
$ws.Range("L127").Value = 5239.1538
$ws.Range("M127").Value = -1862.428400000001
$ws.Range("N127").Value = -15159.1538
$ws.Range("H129").Value = 2713.8635
$ws.Range("I129").Value = 2161.111
$ws.Range("J129").Value = 3096.5386
$ws.Range("K129").Value = 6483.333
$ws.Range("L129").Value = 9289.6158
$ws.Range("M129").Value = -1483.333
$ws.Range("N129").Value = -19289.6158
$ws.Range("H131").Value = 11897
$ws.Range("H134").Value = 31926.295
$ws.Range("J134").Value = 31926.295
$ws.Range("L134").Value = 31926.295
$ws.Range("N134").Value = -42066.295
$ws.Range("H138").Value = 3673.5054
$ws.Range("J138").Value = 3201.9866
$ws.Range("L138").Value = 9605.959800000001
$ws.Range("N138").Value = -19885.9598
$ws.Range("H141").Value = 8050.6816
$ws.Range("I141").Value = 5231.3335
$ws.Range("K141").Value = 15694.0005
$ws.Range("M141").Value = -10514.0005
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6018.8667
$ws.Range("I61").Value = 2396.8333
$ws.Range("K61").Value = 2396.8333
$ws.Range("M61").Value = -2184.8333
$ws.Range("H63").Value = 2199
$ws.Range("I63").Value = 1699
$ws.Range("K63").Value = 1699
$ws.Range("M63").Value = -1013
$ws.Range("H66").Value = 2199
$ws.Range("I66").Value = 1699
$ws.Range("K66").Value = 8495
$ws.Range("M66").Value = -5063
$ws.Range("H110").Value = 836.5909
$ws.Range("I110").Value = 766.61536
$ws.Range("J110").Value = 937.6667
$ws.Range("K110").Value = 766.61536
$ws.Range("L110").Value = 937.6667
$ws.Range("M110").Value = 1278.38464
$ws.Range("N110").Value = -5027.6667
$ws.Range("H122").Value = 8763.791999999999
$ws.Range("I122").Value = 11566.857
$ws.Range("J122").Value = 4839.5
$ws.Range("K122").Value = 34700.571
$ws.Range("L122").Value = 14518.5
$ws.Range("M122").Value = -32250.571
$ws.Range("N122").Value = -19418.5
$ws.Range("H136").Value = 6018.8667
$ws.Range("I136").Value = 2396.8333
$ws.Range("K136").Value = 7190.499899999999
$ws.Range("M136").Value = -4640.499899999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4328.2104
$ws.Range("I105").Value = 3681.111
$ws.Range("K105").Value = 3681.111
$ws.Range("M105").Value = -1934.111
$ws.Range("H132").Value = 99999
$ws.Range("J132").Value = 99999
$ws.Range("L132").Value = 99999
$ws.Range("N132").Value = -110119
$ws.Range("H133").Value = 66242
$ws.Range("J133").Value = 66242
$ws.Range("L133").Value = 66242
$ws.Range("N133").Value = -76362
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 95355
$ws.Range("I31").Value = 2624.25
$ws.Range("K31").Value = 2624.25
$ws.Range("M31").Value = -2329.25
$ws.Range("H34").Value = 95355
$ws.Range("I34").Value = 2624.25
$ws.Range("K34").Value = 2624.25
$ws.Range("M34").Value = -2422.25
$ws.Range("H58").Value = 3028.9333
$ws.Range("I58").Value = 1812.2727
$ws.Range("K58").Value = 1812.2727
$ws.Range("M58").Value = -1609.2727
$ws.Range("H105").Value = 1835.5
$ws.Range("I105").Value = 1895.2222
$ws.Range("K105").Value = 1895.2222
$ws.Range("M105").Value = -148.2221999999999
$ws.Range("H122").Value = 5940.636
$ws.Range("I122").Value = 6837.1665
$ws.Range("J122").Value = 4864.8
$ws.Range("K122").Value = 20511.4995
$ws.Range("L122").Value = 14594.4
$ws.Range("M122").Value = -18061.4995
$ws.Range("N122").Value = -19494.4
$ws.Range("H136").Value = 3028.9333
$ws.Range("I136").Value = 1812.2727
$ws.Range("K136").Value = 5436.8181
$ws.Range("M136").Value = -2886.8181
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 9227.739
$ws.Range("I139").Value = 3146.125
$ws.Range("K139").Value = 9438.375
$ws.Range("M139").Value = -4298.375
$ws.Range("H140").Value = 6583215
$ws.Range("I140").Value = 7354345.5
$ws.Range("K140").Value = 22063036.5
$ws.Range("M140").Value = -22057856.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 444.1579
$ws.Range("I2").Value = 49
$ws.Range("K2").Value = 49
$ws.Range("M2").Value = 64
$ws.Range("H102").Value = 18425.688
$ws.Range("I102").Value = 19254.066
$ws.Range("K102").Value = 19254.066
$ws.Range("M102").Value = -17632.066
$ws.Range("H113").Value = 2233.923
$ws.Range("I113").Value = 2604.7
$ws.Range("J113").Value = 998
$ws.Range("K113").Value = 2604.7
$ws.Range("L113").Value = 998
$ws.Range("M113").Value = -434.6999999999998
$ws.Range("N113").Value = -5338
$ws.Range("H122").Value = 71052.8
$ws.Range("I122").Value = 145296.86
$ws.Range("J122").Value = 6089.25
$ws.Range("K122").Value = 435890.58
$ws.Range("L122").Value = 18267.75
$ws.Range("M122").Value = -433440.58
$ws.Range("N122").Value = -23167.75
$ws.Range("H125").Value = 20000
$ws.Range("J125").Value = 20000
$ws.Range("L125").Value = 20000
$ws.Range("N125").Value = -24920
$ws.Range("H126").Value = 9807.571
$ws.Range("I126").Value = 7667.25
$ws.Range("J126").Value = 12661.333
$ws.Range("K126").Value = 23001.75
$ws.Range("L126").Value = 37983.999
$ws.Range("M126").Value = -20531.75
$ws.Range("N126").Value = -42923.999
$ws.Range("H135").Value = 94965.71000000001
$ws.Range("J135").Value = 94965.71000000001
$ws.Range("L135").Value = 94965.71000000001
$ws.Range("N135").Value = -105105.71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6204
$ws.Range("I7").Value = 6204
$ws.Range("K7").Value = 6204
$ws.Range("M7").Value = -6092
$ws.Range("H40").Value = 2796.9
$ws.Range("I40").Value = 2446
$ws.Range("K40").Value = 2446
$ws.Range("M40").Value = -2310
$ws.Range("H46").Value = 1815.8334
$ws.Range("J46").Value = 2313.2856
$ws.Range("L46").Value = 2313.2856
$ws.Range("N46").Value = -2689.2856
$ws.Range("H68").Value = 2539.4167
$ws.Range("I68").Value = 2547.3
$ws.Range("K68").Value = 2547.3
$ws.Range("M68").Value = -1798.3
$ws.Range("H71").Value = 2539.4167
$ws.Range("I71").Value = 2547.3
$ws.Range("K71").Value = 12736.5
$ws.Range("M71").Value = -8992.5
$ws.Range("H122").Value = 3614.4285
$ws.Range("I122").Value = 3284.9333
$ws.Range("J122").Value = 4438.1665
$ws.Range("K122").Value = 9854.7999
$ws.Range("L122").Value = 13314.4995
$ws.Range("M122").Value = -7404.7999
$ws.Range("N122").Value = -18214.4995
$ws.Range("H126").Value = 6204
$ws.Range("I126").Value = 6204
$ws.Range("K126").Value = 18612
$ws.Range("M126").Value = -16142
$ws.Range("H127").Value = 97904.336
$ws.Range("J127").Value = 97904.336
$ws.Range("L127").Value = 97904.336
$ws.Range("N127").Value = -107824.336
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5870
$ws.Range("I81").Value = 5968.421
$ws.Range("K81").Value = 11936.842
$ws.Range("M81").Value = -10875.842
$ws.Range("H84").Value = 5870
$ws.Range("I84").Value = 5968.421
$ws.Range("K84").Value = 59684.21000000001
$ws.Range("M84").Value = -54380.21000000001
$ws.Range("H122").Value = 89631
$ws.Range("I122").Value = 105428.45
$ws.Range("K122").Value = 316285.35
$ws.Range("M122").Value = -313835.35
$ws.Range("H126").Value = 2287.1765
$ws.Range("I126").Value = 1979.6774
$ws.Range("K126").Value = 5939.0322
$ws.Range("M126").Value = -3469.0322
$ws.Range("H136").Value = 6918.447
$ws.Range("I136").Value = 5655.8887
$ws.Range("J136").Value = 11050.454
$ws.Range("K136").Value = 16967.6661
$ws.Range("L136").Value = 33151.362
$ws.Range("M136").Value = -14417.6661
$ws.Range("N136").Value = -38251.362
